$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.973.43"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.361.08"
$ws.Range("E3").Value = "  +1.64%  "

$ws.Range("E4").Value = "  -0.03%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "303.09"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "95.43"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("E7").Value = "  -0.01%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.502"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.23%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.476"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -3.61%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "34.40"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "

$ws.Range("E11").Value = "  +2.39%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0786"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("E13").Value = "  -1.61%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.70"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.726.18"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").Value = "2.353.67"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "42.976.99"
$ws.Range("E18").Value = "  +0.55%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.92"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("E20").Value = "  +1.45%  "

$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  -0.10%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "67.92"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "235.20"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("E24").Value = "  -1.84%  "

$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  +0.04%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "24.38"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.36"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +15.04%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.32"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.50%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "32.37"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +3.05%  "

$ws.Range("E31").Value = "  -0.12%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.01"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "17.51"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0723"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +3.83%  "

$ws.Range("E35").Value = "  +6.36%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "127.98"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -8.52%  "

$ws.Range("E37").Value = "  +0.60%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "4.33"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("E39").Value = "  +3.73%  "

$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("E41").Value = "  -0.73%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "20.65"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -7.76%  "

$ws.Range("D43").Value = "1.927.11"
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  +3.19%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "9.24"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -9.74%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").Value = "2.587.03"
$ws.Range("E48").Value = "  +1.36%  "

$ws.Range("E49").Value = "  +2.98%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "71.41"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("E51").Value = "  +0.71%  "
